# Insert one new data row above current row 875 (shifting existing rows down
# by one), then populate the newly inserted row with its values.
#
# Resulting layout: a brand-new row "2026/02/26 木 2 201" becomes row 875,
# and everything that used to be rows 875-916 becomes rows 876-917.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$insertRow = 875

# Shift rows 875:916 down to 876:917 by inserting a blank row above 875.
$ws.Rows.Item($insertRow).Insert()

# Populate the newly-inserted row 875 with the new record. The date column
# holds plain text (not a real date), so force text formatting before
# assigning the value to keep Excel from auto-converting it to a date, then
# restore the default "Normal" style so the cell matches its neighbors
# (which all use the default, unstyled format with no explicit style).
$dateCell = $ws.Cells.Item($insertRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/02/26"
$dateCell.Style = "Normal"
$ws.Cells.Item($insertRow, 2).Value = "木"
$ws.Cells.Item($insertRow, 3).Value = 2
$ws.Cells.Item($insertRow, 4).Value = 201
